$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getpages_with_noindex")

$ws.Range("A3").Value = "/about-cancer/coping/self-image"
$ws.Range("B3").Value = "article"
$ws.Range("C3").Value = "English"

$ws.Activate()
$ws.Range("C3").Select()
